$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at 34 (old row 34 "totals" -> 35, old row 35 "footer" -> 36)
$ws.Rows.Item(34).Insert()

# 2. Duplicate row 33 (currently the last item row, "محلول ملح") down into the
#    newly inserted row 34, carrying formatting, values and merged-cell layout.
$ws.Range("A33:Q33").Copy($ws.Range("A34:Q34"))

# 3. Fix up row heights:
#    - new row 34 should be 25.5 (standard data row height)
#    - the totals row (now 35) shrinks from 26.25 to 25.5
$ws.Rows.Item(34).RowHeight = 25.5
$ws.Rows.Item(35).RowHeight = 25.5

# 4. Row 34 currently holds a duplicate of the old row 33 ("محلول ملح" item);
#    give it the correct sequence number (28) to follow row 33 (27).
$ws.Range("A34").Value = 28

# 5. Turn row 33 into the new item "كحول طبي" inserted ahead of "محلول ملح".
$ws.Range("C33").Value = "كحول طبي"
$ws.Range("H33").Value = "12:0"
$ws.Range("N33").Value = "15.00"
$ws.Range("P33").Value = "15.0000"
$ws.Range("Q33").Value = "1:0"

# 6. Update the grand-total cell (now row 35) to include the new item's price.
$ws.Range("N35").Value = 1403.495
